# Update mods data [2025-12-14 15:09:40]
# Appends a new daily row (row 35) to the ModCounts sheet:
#   Date = 2025/12/14, Game = 逃离鸭科夫, ModCount = 1349
# matching the existing table's formatting (center/center alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the formatting of the last existing data row (34) and apply it to
# the new row (35) first, so the appended cells inherit the same style
# (centered horizontally/vertically) as every other data row.
$ws.Range("A34:C34").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's values. The date column stores its values as plain
# text in this workbook (e.g. "2025/12/13"), so prefix with an apostrophe
# to stop Excel from auto-converting the text into a date serial number.
$ws.Range("A35").Value = "'2025/12/14"
$ws.Range("B35").Value = "逃离鸭科夫"
$ws.Range("C35").Value = 1349
